# "fixed bug for swap scene"
# Column K (CanClone) should be 1 for every data row from 11 through 35.
# Row 12 additionally had its Share(J)/CanClone(K) values swapped:
#   J12: 1 -> 0, K12: 0 -> 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the swapped values on row 12 first.
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1

# Set CanClone = 1 for rows 11, 13-35 (row 12 already set above).
$ws.Range("K11").Value = 1
$ws.Range("K13:K35").Value = 1

# Reflect the user's on-screen selection state after the edit.
$ws.Range("K13").Select()
